$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the entire "Legislature" row (row 11) - shifts all following rows up by one.
$ws.Rows.Item(11).Delete()

# "Drinking Water Division (SWRCB)" row (now row 17): B value 2 -> 1
$ws.Range("B17").Value = 1

# "Local Water Boards" row (now row 22): B value 2 -> 1
$ws.Range("B22").Value = 1

# "CV SALTS management zones" row (now row 23): fill C:F with -1,
# matching the style already used on similar rows (e.g. C20, the NRCS row).
$ws.Range("C20").Copy()
$ws.Range("C23:F23").PasteSpecial(-4122)
$ws.Range("C23:F23").Value = -1
$excel.CutCopyMode = $false

# Match the saved selection/active-range state.
$ws.Range("C23:F23").Select() | Out-Null
